$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-28 Saturday" "2026-03-01 Sunday"

Replace-Text "913×2=1826" "441×4=1764"
Replace-Text "661×3=1983" "380×2=760"
Replace-Text "703×3=2109" "818×5=4090"
Replace-Text "358×5=1790" "300×6=1800"
Replace-Text "726×9=6534" "984×9=8856"

Replace-Text "124×6=744" "442×7=3094"
Replace-Text "554×8=4432" "472×5=2360"
Replace-Text "683×9=6147" "217×5=1085"
Replace-Text "506×9=4554" "879×5=4395"
Replace-Text "785×9=7065" "969×5=4845"

Replace-Text "803×8=6424" "349×7=2443"
Replace-Text "180×7=1260" "842×5=4210"
Replace-Text "231×3=693" "176×7=1232"
Replace-Text "152×2=304" "846×2=1692"
Replace-Text "946×6=5676" "823×7=5761"

Replace-Text "392×6=2352" "509×6=3054"
Replace-Text "636×9=5724" "546×6=3276"
Replace-Text "611×7=4277" "179×9=1611"
Replace-Text "659×4=2636" "275×5=1375"
Replace-Text "667×4=2668" "780×8=6240"

Replace-Text "181×7=1267" "470×8=3760"
Replace-Text "961×4=3844" "811×7=5677"
Replace-Text "650×4=2600" "240×8=1920"
Replace-Text "328×9=2952" "383×5=1915"
Replace-Text "592×2=1184" "947×3=2841"
